$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new bullet after "- Served as the technical lead for a
#    greenfield microservice project":
#    "- Adhered to data driven design principles to create a robust,
#     reusable income calculator"
# ---------------------------------------------------------------------
$pLead = $d.Paragraphs(22)
$rLead = $pLead.Range
$rLead.Collapse(0)
$rLead.InsertParagraphAfter()
$pAdhered = $d.Paragraphs(23)
$pAdhered.Range.Text = "- Adhered to data driven design principles to create a robust, reusable income calculator"

# ---------------------------------------------------------------------
# 2. "- Embraced agile principles to deliver incremental value" ->
#    "- Embraced agile methodology to deliver incremental value"
# ---------------------------------------------------------------------
$pAgile = $d.Paragraphs(24)
$pAgile.Range.Find.Execute("principles", $true, $false, $false, $false, $false, $true, 1, $false, "methodology", 2)

# ---------------------------------------------------------------------
# 3. Insert a new bullet after "- Recognized as SQL data abstraction
#    expert for legacy applications":
#    "- Created stored procedures, tables, and question sets to
#     support new pages"
# ---------------------------------------------------------------------
$pRecognized = $d.Paragraphs(25)
$rRecognized = $pRecognized.Range
$rRecognized.Collapse(0)
$rRecognized.InsertParagraphAfter()
$pCreated = $d.Paragraphs(26)
$pCreated.Range.Text = "- Created stored procedures, tables, and question sets to support new pages"

# ---------------------------------------------------------------------
# 4. Relocate the hidden "_GoBack" bookmark (tracks Word's last-edit
#    position) from just after "Quality" to right after the text just
#    inserted above, collapsed (start == end), before the paragraph
#    mark - mirroring where Word itself would leave it after the most
#    recent edit.
#
#    A zero-width range placed exactly at "end of paragraph text,
#    before the pilcrow" cannot be fed to Bookmarks.Add directly in
#    this runtime, so we park a one-character placeholder there first,
#    anchor the bookmark just before it (now a normal mid-paragraph
#    position), and then remove the placeholder.
# ---------------------------------------------------------------------
$pFinal = $d.Paragraphs(26)
$rFinal = $pFinal.Range
$placeholder = $d.Range($rFinal.End - 1, $rFinal.End - 1)
$placeholder.InsertBefore("Z")
$bmPos = $placeholder.Start
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($bmPos, $bmPos + 1).Delete()
